$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2387126.8
$ws.Cells.Item(17, 10).Value = 2387126.8
$ws.Cells.Item(17, 12).Value = 7161380.399999999
$ws.Cells.Item(17, 14).Value = -7161716.399999999

$ws.Cells.Item(76, 8).Value = 3679.9678

$ws.Cells.Item(79, 8).Value = 3679.9678

$ws.Cells.Item(112, 8).Value = 1209.3214
$ws.Cells.Item(112, 9).Value = 900
$ws.Cells.Item(112, 10).Value = 1233.1154
$ws.Cells.Item(112, 11).Value = 2700
$ws.Cells.Item(112, 12).Value = 3699.3462
$ws.Cells.Item(112, 13).Value = -1592
$ws.Cells.Item(112, 14).Value = -5915.3462

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1844.7925
$ws.Cells.Item(132, 9).Value = 1474.2703
$ws.Cells.Item(132, 11).Value = 4422.810899999999
$ws.Cells.Item(132, 13).Value = -1892.810899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 3643.3333
$ws.Cells.Item(105, 9).Value = 4000
$ws.Cells.Item(105, 11).Value = 4000
$ws.Cells.Item(105, 13).Value = -2253

$ws.Cells.Item(141, 8).Value = 44162.727
$ws.Cells.Item(141, 10).Value = 44162.727
$ws.Cells.Item(141, 12).Value = 44162.727
$ws.Cells.Item(141, 14).Value = -54522.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1528.1111
$ws.Cells.Item(16, 9).Value = 1567.5
$ws.Cells.Item(16, 10).Value = 1213
$ws.Cells.Item(16, 11).Value = 1567.5
$ws.Cells.Item(16, 12).Value = 1213
$ws.Cells.Item(16, 13).Value = -1280.5
$ws.Cells.Item(16, 14).Value = -1787

$ws.Cells.Item(113, 8).Value = 1528.1111
$ws.Cells.Item(113, 9).Value = 1567.5
$ws.Cells.Item(113, 10).Value = 1213
$ws.Cells.Item(113, 11).Value = 1567.5
$ws.Cells.Item(113, 12).Value = 1213
$ws.Cells.Item(113, 13).Value = 602.5
$ws.Cells.Item(113, 14).Value = -5553

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 9).Value = 440.66666
$ws.Cells.Item(5, 11).Value = 1321.99998
$ws.Cells.Item(5, 13).Value = -1209.99998

$ws.Cells.Item(39, 8).Value = 11801.607
$ws.Cells.Item(39, 9).Value = 52451.5
$ws.Cells.Item(39, 10).Value = 10296.056
$ws.Cells.Item(39, 11).Value = 157354.5
$ws.Cells.Item(39, 12).Value = 30888.168
$ws.Cells.Item(39, 13).Value = -157060.5
$ws.Cells.Item(39, 14).Value = -31476.168

$ws.Cells.Item(81, 8).Value = 5276.6665
$ws.Cells.Item(81, 9).Value = 800
$ws.Cells.Item(81, 10).Value = 7515
$ws.Cells.Item(81, 11).Value = 2400
$ws.Cells.Item(81, 12).Value = 22545
$ws.Cells.Item(81, 13).Value = -1277
$ws.Cells.Item(81, 14).Value = -24791

$ws.Cells.Item(84, 8).Value = 5276.6665
$ws.Cells.Item(84, 9).Value = 800
$ws.Cells.Item(84, 10).Value = 7515
$ws.Cells.Item(84, 11).Value = 7200
$ws.Cells.Item(84, 12).Value = 67635
$ws.Cells.Item(84, 13).Value = -1584
$ws.Cells.Item(84, 14).Value = -78867

$ws.Cells.Item(114, 8).Value = 253.18182
$ws.Cells.Item(114, 9).Value = 225.4
$ws.Cells.Item(114, 10).Value = 531
$ws.Cells.Item(114, 11).Value = 676.2
$ws.Cells.Item(114, 12).Value = 1593
$ws.Cells.Item(114, 13).Value = 2577.8
$ws.Cells.Item(114, 14).Value = -8101

$ws.Cells.Item(129, 8).Value = 1470.5294
$ws.Cells.Item(129, 9).Value = 673.26666
$ws.Cells.Item(129, 10).Value = 2099.9473
$ws.Cells.Item(129, 11).Value = 2019.79998
$ws.Cells.Item(129, 12).Value = 6299.841899999999
$ws.Cells.Item(129, 13).Value = 2980.20002
$ws.Cells.Item(129, 14).Value = -16299.8419

$ws.Cells.Item(131, 8).Value = 859.48
$ws.Cells.Item(131, 9).Value = 295
$ws.Cells.Item(131, 10).Value = 883
$ws.Cells.Item(131, 11).Value = 885
$ws.Cells.Item(131, 12).Value = 2649
$ws.Cells.Item(131, 13).Value = 4155
$ws.Cells.Item(131, 14).Value = -12729

$ws.Cells.Item(135, 9).Value = 440.66666
$ws.Cells.Item(135, 11).Value = 3965.99994
$ws.Cells.Item(135, 13).Value = -1430.99994

$ws.Cells.Item(137, 8).Value = 2004.5714
$ws.Cells.Item(137, 10).Value = 2077.7144
$ws.Cells.Item(137, 12).Value = 6233.1432
$ws.Cells.Item(137, 14).Value = -16433.1432

$ws.Cells.Item(140, 8).Value = 1945
$ws.Cells.Item(140, 9).Value = 1890
$ws.Cells.Item(140, 10).Value = 2000
$ws.Cells.Item(140, 11).Value = 5670
$ws.Cells.Item(140, 12).Value = 6000
$ws.Cells.Item(140, 13).Value = -490
$ws.Cells.Item(140, 14).Value = -16360

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value = 205125
$ws.Cells.Item(21, 9).Value = 6800
$ws.Cells.Item(21, 10).Value = 535666.7
$ws.Cells.Item(21, 11).Value = 6800
$ws.Cells.Item(21, 12).Value = 535666.7
$ws.Cells.Item(21, 13).Value = -6627
$ws.Cells.Item(21, 14).Value = -536012.7

$ws.Cells.Item(30, 8).Value = 205125
$ws.Cells.Item(30, 9).Value = 6800
$ws.Cells.Item(30, 10).Value = 535666.7
$ws.Cells.Item(30, 11).Value = 6800
$ws.Cells.Item(30, 12).Value = 535666.7
$ws.Cells.Item(30, 13).Value = -6695
$ws.Cells.Item(30, 14).Value = -535876.7

$ws.Cells.Item(70, 8).Value = 5048.62
$ws.Cells.Item(70, 9).Value = 4862.6113
$ws.Cells.Item(70, 10).Value = 5153.25
$ws.Cells.Item(70, 11).Value = 4862.6113
$ws.Cells.Item(70, 12).Value = 5153.25
$ws.Cells.Item(70, 13).Value = -4592.6113
$ws.Cells.Item(70, 14).Value = -5693.25

$ws.Cells.Item(73, 8).Value = 5048.62
$ws.Cells.Item(73, 9).Value = 4862.6113
$ws.Cells.Item(73, 10).Value = 5153.25
$ws.Cells.Item(73, 11).Value = 4862.6113
$ws.Cells.Item(73, 12).Value = 5153.25
$ws.Cells.Item(73, 13).Value = -3926.6113
$ws.Cells.Item(73, 14).Value = -7025.25

$ws.Cells.Item(80, 8).Value = 2728.4614
$ws.Cells.Item(80, 9).Value = 2770
$ws.Cells.Item(80, 11).Value = 2770
$ws.Cells.Item(80, 13).Value = -1772

$ws.Cells.Item(83, 8).Value = 2728.4614
$ws.Cells.Item(83, 9).Value = 2770
$ws.Cells.Item(83, 11).Value = 13850
$ws.Cells.Item(83, 13).Value = -8858

$ws.Cells.Item(123, 8).Value = 18386.223
$ws.Cells.Item(123, 10).Value = 18386.223
$ws.Cells.Item(123, 12).Value = 18386.223
$ws.Cells.Item(123, 14).Value = -23286.223

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(23, 8).Value = 10533.333
$ws.Cells.Item(23, 9).Value = 5800
$ws.Cells.Item(23, 11).Value = 5800
$ws.Cells.Item(23, 13).Value = -5570

$ws.Cells.Item(34, 8).Value = 50000
$ws.Cells.Item(34, 10).Value = 50000
$ws.Cells.Item(34, 12).Value = 50000
$ws.Cells.Item(34, 14).Value = -50344

$ws.Cells.Item(61, 8).Value = 1954
$ws.Cells.Item(61, 9).Value = 2064
$ws.Cells.Item(61, 10).Value = 1800
$ws.Cells.Item(61, 11).Value = 2064
$ws.Cells.Item(61, 12).Value = 1800
$ws.Cells.Item(61, 13).Value = -1862
$ws.Cells.Item(61, 14).Value = -2204

$ws.Cells.Item(113, 8).Value = 1954
$ws.Cells.Item(113, 9).Value = 2064
$ws.Cells.Item(113, 10).Value = 1800
$ws.Cells.Item(113, 11).Value = 2064
$ws.Cells.Item(113, 12).Value = 1800
$ws.Cells.Item(113, 13).Value = 106
$ws.Cells.Item(113, 14).Value = -6140
